$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.452.31"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "3.514.24"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'600.97"
$ws.Range("D6").Value = "'176.67"
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "'7.15"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").Value = "'0.433"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "4.120.31"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'31.45"
$ws.Range("E13").Value = "  +11.11%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "67.465.82"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "3.502.84"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "'6.33"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'14.70"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").Value = "'394.95"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'73.59"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'0.541"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D25").Value = "'5.70"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'10.30"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'23.78"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "'163.78"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'7.07"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("D40").Value = "'4.69"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'26.61"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'27.32"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0736"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "2.809.28"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").Value = "'42.58"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").Value = "'338.87"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").Value = "'34.03"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  -0.11%  "
